$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.798.77'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.628.38'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.40'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5041'
$ws.Range("E6").Value = '  -0.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2562'
$ws.Range("E8").Value = '  +0.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06323'
$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07751'
$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.224'
$ws.Range("E12").Value = '  -1.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.627.03'
$ws.Range("E13").Value = '  -0.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5466'
$ws.Range("E14").Value = '  +0.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.57'
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7631'
$ws.Range("E16").Value = '  -2.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.828.80'
$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.415'
$ws.Range("E19").Value = '  -0.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.93'
$ws.Range("E20").Value = '  -2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.867'
$ws.Range("E21").Value = '  -0.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.026'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  +2.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.94'
$ws.Range("E25").Value = '  +0.86%  '

$ws.Range("E26").Value = '  +4.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.772'
$ws.Range("E27").Value = '  -1.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.48'
$ws.Range("E28").Value = '  -1.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.236'
$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04857'
$ws.Range("E30").Value = '  -1.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.232'
$ws.Range("E31").Value = '  -0.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.174'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.535'
$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.371'
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("E36").Value = '  -1.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5471'
$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.117.18'
$ws.Range("E38").Value = '  -1.29%  '

$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.553'
$ws.Range("E41").Value = '  -0.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7956'
$ws.Range("E42").Value = '  -2.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.85'
$ws.Range("E43").Value = '  -2.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₈118'
$ws.Range("E44").Value = '  -7.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.761.97'
$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4440'
$ws.Range("E46").Value = '  -2.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.56'
$ws.Range("E48").Value = '  -0.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05128'
$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.525'
$ws.Range("E50").Value = '  +2.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  -0.15%  '

